$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $cell = $Sheet.Range($Addr)
    $origFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.NumberFormat = $origFormat
}

Set-TextValue $ws "D2" "26.588.27"
Set-TextValue $ws "E2" "  -2.13%  "
Set-TextValue $ws "D3" "1.583.60"
Set-TextValue $ws "E3" "  -2.85%  "
Set-TextValue $ws "E4" "  -0.06%  "
Set-TextValue $ws "D5" "210.51"
Set-TextValue $ws "E5" "  -2.52%  "
Set-TextValue $ws "D6" "0.507"
Set-TextValue $ws "E6" "  -2.07%  "
Set-TextValue $ws "E7" "  -0.05%  "
Set-TextValue $ws "E8" "  -2.55%  "
Set-TextValue $ws "E9" "  -1.02%  "
Set-TextValue $ws "D10" "19.54"
Set-TextValue $ws "E10" "  -3.55%  "
Set-TextValue $ws "E11" "  -1.76%  "
Set-TextValue $ws "D12" "1.805.87"
Set-TextValue $ws "E12" "  -2.81%  "
Set-TextValue $ws "D13" "1.578.15"
Set-TextValue $ws "E13" "  -3.18%  "
Set-TextValue $ws "E14" "  -1.78%  "
Set-TextValue $ws "D15" "0.527"
Set-TextValue $ws "E15" "  -2.80%  "
Set-TextValue $ws "D16" "64.54"
Set-TextValue $ws "E16" "  +0.00%  "
Set-TextValue $ws "D17" "26.609.51"
Set-TextValue $ws "E17" "  -2.04%  "
Set-TextValue $ws "D18" "0.0₃0729"
Set-TextValue $ws "E18" "  -0.44%  "
Set-TextValue $ws "B19" "Dai"
Set-TextValue $ws "C19" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D19" "1.00"
Set-TextValue $ws "E19" "  +0.04%  "
Set-TextValue $ws "B20" "BitcoinCash"
Set-TextValue $ws "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D20" "207.64"
Set-TextValue $ws "E20" "  -3.75%  "
Set-TextValue $ws "E21" "  -1.78%  "
Set-TextValue $ws "E22" "  -3.23%  "
Set-TextValue $ws "E23" "  -4.14%  "
Set-TextValue $ws "E24" "  -2.37%  "
Set-TextValue $ws "D25" "146.42"
Set-TextValue $ws "E25" "  -1.02%  "
Set-TextValue $ws "B26" "BinanceUSD"
Set-TextValue $ws "C26" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws "D26" "1.00"
Set-TextValue $ws "E26" "  -0.09%  "
Set-TextValue $ws "B27" "Cosmos"
Set-TextValue $ws "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D27" "7.40"
Set-TextValue $ws "E27" "  +1.78%  "
Set-TextValue $ws "E28" "  -4.35%  "
Set-TextValue $ws "D29" "15.28"
Set-TextValue $ws "E29" "  -1.86%  "
Set-TextValue $ws "D30" "0.0503"
Set-TextValue $ws "E30" "  -0.49%  "
Set-TextValue $ws "E31" "  -2.11%  "
Set-TextValue $ws "D32" "3.25"
Set-TextValue $ws "E32" "  -4.10%  "
Set-TextValue $ws "D33" "0.682"
Set-TextValue $ws "E33" "  +25.96%  "
Set-TextValue $ws "D35" "1.320.38"
Set-TextValue $ws "E35" "  +0.62%  "
Set-TextValue $ws "D36" "2.49"
Set-TextValue $ws "E36" "  +1.22%  "
Set-TextValue $ws "D37" "1.50"
Set-TextValue $ws "E37" "  -3.73%  "
Set-TextValue $ws "E38" "  -1.29%  "
Set-TextValue $ws "D39" "0.823"
Set-TextValue $ws "E39" "  -3.18%  "
Set-TextValue $ws "E40" "  -0.04%  "
Set-TextValue $ws "E41" "  +2.80%  "
Set-TextValue $ws "D42" "0.783"
Set-TextValue $ws "E42" "  -2.25%  "
Set-TextValue $ws "E43" "  -3.58%  "
Set-TextValue $ws "D44" "63.44"
Set-TextValue $ws "E44" "  -0.43%  "
Set-TextValue $ws "D45" "1.718.65"
Set-TextValue $ws "E45" "  -2.72%  "
Set-TextValue $ws "D46" "89.46"
Set-TextValue $ws "E46" "  -1.37%  "
Set-TextValue $ws "E47" "  +1.11%  "
Set-TextValue $ws "D48" "0.830"
Set-TextValue $ws "E48" "  +3.41%  "
Set-TextValue $ws "D49" "0.0989"
Set-TextValue $ws "E49" "  +3.49%  "
Set-TextValue $ws "D50" "0.0507"
Set-TextValue $ws "E50" "  -1.76%  "
Set-TextValue $ws "E51" "  -0.90%  "
